$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores plain looking numeric text (e.g. "1.003") as
# an inline string in the source workbook. Mark it as Text first so Excel
# does not silently convert these values into real numbers when the .Value
# property is set below.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '27.867.02'
$ws.Cells.Item(2, 5).Value = '  +0.46%  '
$ws.Cells.Item(3, 4).Value = '1.769.06'
$ws.Cells.Item(3, 5).Value = '  +0.61%  '
$ws.Cells.Item(4, 4).Value = '1.003'
$ws.Cells.Item(4, 5).Value = '  +0.33%  '
$ws.Cells.Item(5, 4).Value = '327.93'
$ws.Cells.Item(5, 5).Value = '  +0.89%  '
$ws.Cells.Item(6, 4).Value = '1.002'
$ws.Cells.Item(6, 5).Value = '  +0.40%  '
$ws.Cells.Item(7, 4).Value = '0.4482'
$ws.Cells.Item(7, 5).Value = '  -2.40%  '
$ws.Cells.Item(8, 4).Value = '0.3571'
$ws.Cells.Item(8, 5).Value = '  -0.92%  '
$ws.Cells.Item(9, 4).Value = '42.14'
$ws.Cells.Item(9, 5).Value = '  -0.11%  '
$ws.Cells.Item(10, 4).Value = '0.07435'
$ws.Cells.Item(10, 5).Value = '  -1.35%  '
$ws.Cells.Item(11, 4).Value = '1.094'
$ws.Cells.Item(11, 5).Value = '  -0.75%  '
$ws.Cells.Item(12, 4).Value = '1.002'
$ws.Cells.Item(12, 5).Value = '  +0.33%  '
$ws.Cells.Item(13, 4).Value = '20.88'
$ws.Cells.Item(13, 5).Value = '  -0.26%  '
$ws.Cells.Item(14, 4).Value = '6.049'
$ws.Cells.Item(14, 5).Value = '  +0.32%  '
$ws.Cells.Item(15, 4).Value = '7.216'
$ws.Cells.Item(15, 5).Value = '  +1.11%  '
$ws.Cells.Item(16, 4).Value = '1.781.24'
$ws.Cells.Item(16, 5).Value = '  +1.51%  '
$ws.Cells.Item(17, 4).Value = '92.99'
$ws.Cells.Item(17, 5).Value = '  +0.22%  '
$ws.Cells.Item(18, 4).Value = '0.00001059'
$ws.Cells.Item(18, 5).Value = '  -0.77%  '
$ws.Cells.Item(19, 4).Value = '0.06417'
$ws.Cells.Item(19, 5).Value = '  -0.04%  '
$ws.Cells.Item(20, 5).Value = '  +0.37%  '
$ws.Cells.Item(21, 4).Value = '17.19'
$ws.Cells.Item(21, 5).Value = '  +2.07%  '
$ws.Cells.Item(22, 4).Value = '5.822'
$ws.Cells.Item(22, 5).Value = '  -0.16%  '
$ws.Cells.Item(23, 4).Value = '27.879.72'
$ws.Cells.Item(23, 5).Value = '  +0.38%  '
$ws.Cells.Item(24, 4).Value = '11.32'
$ws.Cells.Item(24, 5).Value = '  +0.49%  '
$ws.Cells.Item(25, 4).Value = '2.118'
$ws.Cells.Item(25, 5).Value = '  +0.49%  '
$ws.Cells.Item(26, 4).Value = '162.56'
$ws.Cells.Item(26, 5).Value = '  -0.59%  '
$ws.Cells.Item(27, 4).Value = '20.23'
$ws.Cells.Item(27, 5).Value = '  -1.21%  '
$ws.Cells.Item(28, 4).Value = '1.980.06'
$ws.Cells.Item(28, 5).Value = '  +1.15%  '
$ws.Cells.Item(29, 4).Value = '2.171'
$ws.Cells.Item(29, 5).Value = '  +3.13%  '
$ws.Cells.Item(30, 4).Value = '125.40'
$ws.Cells.Item(30, 5).Value = '  -1.23%  '
$ws.Cells.Item(31, 4).Value = '1.101'
$ws.Cells.Item(31, 5).Value = '  +1.96%  '
$ws.Cells.Item(32, 4).Value = '0.09133'
$ws.Cells.Item(32, 5).Value = '  -0.82%  '
$ws.Cells.Item(33, 4).Value = '5.588'
$ws.Cells.Item(33, 5).Value = '  +0.42%  '
$ws.Cells.Item(34, 4).Value = '3.633'
$ws.Cells.Item(34, 5).Value = '  -0.87%  '
$ws.Cells.Item(35, 4).Value = '11.83'
$ws.Cells.Item(35, 5).Value = '  -1.23%  '
$ws.Cells.Item(36, 4).Value = '0.02293'
$ws.Cells.Item(36, 5).Value = '  -0.52%  '
$ws.Cells.Item(37, 4).Value = '0.06088'
$ws.Cells.Item(37, 5).Value = '  +0.58%  '
$ws.Cells.Item(38, 4).Value = '0.2098'
$ws.Cells.Item(38, 5).Value = '  -0.35%  '
$ws.Cells.Item(39, 4).Value = '0.6329'
$ws.Cells.Item(39, 5).Value = '  -1.10%  '
$ws.Cells.Item(40, 4).Value = '4.956'
$ws.Cells.Item(40, 5).Value = '  -0.60%  '
$ws.Cells.Item(41, 4).Value = '1.182'
$ws.Cells.Item(41, 5).Value = '  -1.83%  '
$ws.Cells.Item(42, 5).Value = '  +1.16%  '
$ws.Cells.Item(43, 4).Value = '7.926'
$ws.Cells.Item(43, 5).Value = '  +1.17%  '
$ws.Cells.Item(44, 4).Value = '13.21'
$ws.Cells.Item(44, 5).Value = '  -0.82%  '
$ws.Cells.Item(45, 5).Value = '  +0.67%  '
$ws.Cells.Item(46, 4).Value = '0.5873'
$ws.Cells.Item(46, 5).Value = '  -1.05%  '
$ws.Cells.Item(47, 4).Value = '122.60'
$ws.Cells.Item(47, 5).Value = '  -0.52%  '
$ws.Cells.Item(48, 4).Value = '1.957'
$ws.Cells.Item(48, 5).Value = '  -0.40%  '
$ws.Cells.Item(49, 4).Value = '0.06909'
$ws.Cells.Item(49, 5).Value = '  +0.45%  '
$ws.Cells.Item(50, 4).Value = '1.139'
$ws.Cells.Item(50, 5).Value = '  -0.99%  '
$ws.Cells.Item(51, 4).Value = '72.90'
$ws.Cells.Item(51, 5).Value = '  +0.58%  '

# Re-apply the default (unstyled) look to the Price column now that the
# values have been written, so no extra cell styling is introduced.
$ws.Range("D2:D51").Style = "Normal"
